$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 27, shifting current rows 27:49 down to 28:50.
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 (same static fields as its neighbours,
# new date + new price figures).
$ws.Cells.Item(27, 1).Value2 = 1
$ws.Cells.Item(27, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(27, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(27, 4).Value2 = 44586
$ws.Cells.Item(27, 5).Value2 = 15
$ws.Cells.Item(27, 6).Value2 = 100112031
$ws.Cells.Item(27, 7).Value2 = "Poroto verde"
$ws.Cells.Item(27, 8).Value2 = "Sin especificar"
$ws.Cells.Item(27, 9).Value2 = "Primera"
$ws.Cells.Item(27, 10).Value2 = 1500
$ws.Cells.Item(27, 11).Value2 = 1400
$ws.Cells.Item(27, 12).Value2 = 1500
$ws.Cells.Item(27, 13).Value2 = 1450
$ws.Cells.Item(27, 14).Value2 = "`$/kilo"
$ws.Cells.Item(27, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(27, 16).Value2 = 1450
$ws.Cells.Item(27, 17).Value2 = 1
$ws.Cells.Item(27, 18).Value2 = "Hortaliza"
